# Apply the two changes described by the diff:
#  1. DataFile sheet: add "zip" to the file_format dropdown list (column C).
#  2. OntologyTerm sheet: OntologyTerm now inherits id/title/description from
#     NamedThing, so the columns are reordered to
#     label, definition, ontology, id, title, description (A1:F1).

$wb = $excel.ActiveWorkbook

# --- 1. Extend the file_format list validation on the DataFile sheet ---
$dataFileSheet = $wb.Worksheets.Item("DataFile")
$formatValidation = $dataFileSheet.Range("C2:C1048576").Validation
$formatValidation.Formula1 = '"mrc,tiff,hdf5,star,pdb,mmcif,mtz,cbf,ascii,thermo_raw,zip"'

# --- 2. Rebuild the OntologyTerm header row to inherit from NamedThing ---
$ontologyTermSheet = $wb.Worksheets.Item("OntologyTerm")
$ontologyTermSheet.Range("A1").Value = "label"
$ontologyTermSheet.Range("B1").Value = "definition"
$ontologyTermSheet.Range("C1").Value = "ontology"
$ontologyTermSheet.Range("D1").Value = "id"
$ontologyTermSheet.Range("E1").Value = "title"
$ontologyTermSheet.Range("F1").Value = "description"
